# Test Data.xlsx refresh:
#  - Sheet1 gets populated with the full 11-employee data set (same rows
#    already present on Sheet2).
#  - Sheet2's "Extra Duty Allowance" (column J) scratch column is removed.
#  - The "Test Name" sheet gets a freshly generated test name in A2.
#  - Sheet1 becomes the active/selected sheet instead of "Test Name".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Test Name")

# --- Sheet1: bring in the rest of the employee rows (5-12) from Sheet2 ---
$srcRows = $ws2.Range("A5:I12")
$dstRows = $ws1.Range("A5:I12")
$dstRows.Value2 = $srcRows.Value2

# --- Sheet2: drop the now-unused "Extra Duty Allowance" column (J) ---
$ws2.Range("J1:J12").Clear()

# --- Test Name sheet: refresh the generated test employee name ---
$ws3.Range("A2").Value = "MHArryMt"

# --- View state: Sheet1 becomes the active tab/selection ---
$ws1.Activate()
$ws1.Range("A1:I12").Select()
$ws2.Range("A2:I12").Select()
$ws1.Activate()
